$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value2 = 'join course with department on course . deptId equals department . id AND THEN FROM course table WHAT ARE THE count ( title ) , department . deptId AND THEN WHAT IS count ( id ) , course . departmentName , id and _CN_ FROM department AND THEN LIMIT 44 GROUP AUTOMATICALLY'
$ws.Range('B3').Value2 = 'IN Salaries as Ss where Salary is equal to 354 and Salary is greater than 510 WHAT IS THE Salaries . FromDate and EmployeeNumber'
$ws.Range('B4').Value2 = 'IN the staff table WHAT ARE staff . username and count ( password ) and picture GROUP BY AUTOMATIC'
$ws.Range('B5').Value2 = 'FIND count ( endTime ) , the count of ( termId ) FROM ( FIND term . year and count ( endDate ) , term . endDate IN the term table where termPeriod IN ( DISPLAY id FROM TABLE term ) and endDate equals bar GROUP AUTOMATIC ) as sq_alias GROUP AUTOMATIC'
$ws.Range('B6').Value2 = 'WHAT IS THE DepartmentManager . EmployeeNumber , count ( EmployeeNumber ) FROM table DepartmentManager as Dr where EmployeeNumber is equal to foo and EmployeeNumber equals bar GROUP AUTOMATICALLY'
$ws.Range('B7').Value2 = 'natural join city with countrylanguage AND THEN SHOW ME District and the count of ( city_id ) IN city AND THEN SELECT city . IsOfficial , count ( Language ) and countrylanguage . CountryCode FROM table countrylanguage AND THEN HAVING COUNT ( CountryCode ) greater than 112 GROUP BY AUTOMATIC'
$ws.Range('B8').Value2 = 'natural join city with countrylanguage AND THEN IN city as cy WHAT IS THE city . District and ID , the sum of ( Population ) , last_update where last_update equals bar or CountryCode equals bar AND THEN SHOW ME the sum of ( Percentage ) , count ( CountryCode ) and countrylanguage . CountryCode where Language equals bar IN countrylanguage table AND THEN HAVING AVG ( CountryCode ) greater than 191 GROUP AUTOMATIC'
$ws.Range('B9').Value2 = 'IN DepartmentManager as Dr where EmployeeNumber equals bar and FromDate IN ( SELECT DepartmentNumber FROM TABLE DepartmentManager ) PRESENT the count of ( ToDate ) and DepartmentManager . ToDate and count ( FromDate ) GROUP AUTOMATICALLY'
$ws.Range('B10').Value2 = 'natural join customer with store AND THEN IN customer WHAT IS first_name , count ( customer_id ) AND THEN RETRIEVE store . last_update and store_id IN store AND THEN HAVING AVG ( last_update ) equal to 54 GROUP BY AUTOMATICALLY'
$ws.Range('B11').Value2 = 'WHAT ARE THE Title , ToDate IN Titles table where Title is equal to foo or EmployeeNumber is equal to foo HAVING SUM ( Title ) less than 139'
$ws.Range('B12').Value2 = 'IN courseoffering as cg DISPLAY endTime and facultyName HAVING COUNT ( endTime ) equal to 59'
$ws.Range('B13').Value2 = 'WHAT IS term . id , term . endDate FROM term where id is equal to bar and startDate equals bar AND THEN natural join term with term AND THEN natural join term with term AND THEN natural join term with term LIMIT 45 HAVING SUM ( endDate ) greater than 6'
$ws.Range('B14').Value2 = 'SHOW ME term . startDate and the count of ( termPeriod ) and term . year IN ( RETRIEVE count ( endTime ) and termId FROM table courseoffering HAVING COUNT ( termId ) less than 114 GROUP AUTOMATICALLY ) as sq_aliasGROUP BY AUTOMATICALLY'
$ws.Range('B15').Value2 = 'RETRIEVE last_update and city IN the city table HAVING AVG ( city ) less than 43'
$ws.Range('B16').Value2 = 'FROM DepartmentManager as Dr SELECT FromDate and ToDate and the count of ( ToDate ) where FromDate equals foo and ToDate is equal to bar AND THEN DISPLAY count ( DepartmentNumber ) , DepartmentManager . DepartmentNumber IN ( PRESENT the count of ( EmployeeNumber ) , Salaries . Salary where EmployeeNumber equals foo FROM table Salaries as Ss HAVING AVG ( Salary ) less than 38 GROUP AUTOMATIC ) as sq_alias AND THEN natural join DepartmentManager with departments AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B17').Value2 = 'natural join term with term AND THEN FROM term table SHOW ME term . year , term . startDate'
$ws.Range('B18').Value2 = 'IN table film as fm PRESENT count ( description ) and language_id GROUP BY AUTOMATICALLY'
$ws.Range('B19').Value2 = 'WHAT IS THE the count of ( onDays ) , the count of ( courseId ) , roomId FROM the courseoffering table as cg where id equals bar GROUP AUTOMATICALLY'
$ws.Range('B20').Value2 = 'WHAT ARE film . release_year , film . last_update FROM the film table where film_id is equal to bar or last_update equals bar'
$ws.Range('B21').Value2 = 'PRESENT staff . active , the sum of ( active ) IN table staff GROUP AUTOMATIC'
$ws.Range('B22').Value2 = 'IN term table DISPLAY year , term . termPeriod'
$ws.Range('B23').Value2 = 'RETRIEVE country . city_id , country . CountryCode where Name is equal to bar and country_id equals bar FROM city AND THEN WHAT ARE THE average ( GNPOld ) and the average of ( GNP ) and country . Continent FROM country as cy AND THEN natural join city with country AND THEN LIMIT 45 GROUP AUTOMATICALLY'
$ws.Range('B24').Value2 = 'RETRIEVE manager_staff_id , store . last_update FROM the store table as se where last_update equals bar or address_id equals bar'
$ws.Range('B25').Value2 = 'WHAT IS city . city , average ( Population ) where city_id is equal to bar and city_id is equal to bar IN the city table GROUP BY AUTOMATIC'
$ws.Range('B26').Value2 = 'join courseoffering with room on courseoffering . roomId equals room . id AND THEN join room with building on room . buildingId equals building . id AND THEN PRESENT courseId and facultyName , average ( capacity ) where startTime is equal to foo and id equals bar FROM table courseoffering as cg AND THEN IN room DISPLAY building . buildingId , building . wheelchairSpaces AND THEN IN table building WHAT IS count ( buildingNumber ) and the count of ( id ) and id where buildingName IN ( FIND buildingNumber FROM TABLE building ) AND THEN HAVING SUM ( id ) greater than 65 GROUP BY AUTOMATICALLY'
$ws.Range('B27').Value2 = 'natural join city with countrylanguage AND THEN SELECT CountryCode and the sum of ( Population ) FROM table city as cy where country_id equals bar or ID is equal to bar AND THEN IN ( PRESENT city . District and city . city_id and city . city IN table city ) as sq_alias where IsOfficial is equal to bar or Percentage is equal to 1870 FIND countrylanguage . Language , IsOfficial AND THEN ORDER BY Language DESC GROUP BY AUTOMATICALLY'
$ws.Range('B28').Value2 = 'IN table room as rm where roomNumber IN ( FIND buildingId FROM TABLE room ) and floor IN ( FIND area FROM TABLE room ) WHAT ARE id and the sum of ( capacity ) LIMIT 33 GROUP BY AUTOMATICALLY'
$ws.Range('B29').Value2 = 'natural join Salaries with Employees AND THEN FROM table Salaries as Ss SELECT Salaries . EmployeeNumber , the count of ( ToDate ) AND THEN RETRIEVE LastName , EmployeeNumber where Gender is equal to foo or LastName is equal to bar FROM table Employees as Es AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B30').Value2 = 'PRESENT count ( store_id ) and inventory . inventory_id and last_update where store_id equals bar IN the inventory table as iy HAVING AVG ( last_update ) less than 63 GROUP AUTOMATIC'
$ws.Range('B31').Value2 = 'WHAT ARE THE average ( length ) and the count of ( last_update ) and language . title FROM film table AND THEN IN the language table FIND language_id , name where language_id is equal to foo and name equals bar AND THEN natural join film with language AND THEN GROUP BY AUTOMATIC'
$ws.Range('B32').Value2 = 'IN table term as tm where year greater than 1716 WHAT IS THE the sum of ( year ) , term . termPeriod and term . id GROUP BY AUTOMATICALLY LIMIT 43'
$ws.Range('B33').Value2 = 'natural join Salaries with Employees AND THEN IN table Salaries where EmployeeNumber is equal to foo or Salary is less than 1239 WHAT IS THE FromDate and Salary and Salaries . ToDate and EmployeeNumber AND THEN IN Employees table where HireDate is equal to bar WHAT ARE the count of ( FirstName ) and Employees . BirthDate AND THEN GROUP AUTOMATIC LIMIT 12'
$ws.Range('B34').Value2 = 'natural join store with address AND THEN natural join address with city AND THEN FROM the store table GET manager_staff_id and city . store_id AND THEN RETRIEVE postal_code and store . district and phone FROM the address table AND THEN IN ( FROM ( WHAT IS THE city_id , count ( address2 ) FROM table address as as GROUP AUTOMATIC HAVING SUM ( city_id ) equal to 75 ) as sq_alias FIND store . manager_staff_id and count ( manager_staff_id ) GROUP BY AUTOMATIC ) as sq_alias where city is equal to foo or Name is equal to bar FIND city . CountryCode and average ( Population ) AND THEN GROUP BY AUTOMATICALLY HAVING COUNT ( CountryCode ) less than 47'
$ws.Range('B35').Value2 = 'FROM customer table RETRIEVE average ( active ) , the count of ( email ) , count ( last_name ) , active where store_id equals bar and customer_id is equal to bar ORDER BY active DESC GROUP AUTOMATICALLY'
$ws.Range('B36').Value2 = 'FROM Salaries table where FromDate IN ( DISPLAY ToDate IN TABLE Employees ) and Salary equal to 632 WHAT ARE THE the average of ( Salary ) and Salaries . ToDate AND THEN FROM the Employees table RETRIEVE EmployeeNumber and LastName AND THEN natural join Salaries with Employees AND THEN GROUP AUTOMATIC'
$ws.Range('B37').Value2 = 'FROM Salaries table SELECT FromDate and sum ( Salary ) GROUP BY AUTOMATICALLY'
$ws.Range('B38').Value2 = 'natural join Titles with Employees AND THEN PRESENT ToDate and Titles . Title , count ( FromDate ) where ToDate IN ( FIND ToDate IN TABLE Employees ) IN the Titles table AND THEN GET the count of ( HireDate ) and the count of ( HireDate ) IN Employees table AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B39').Value2 = 'DISPLAY sum ( Population ) and city . last_update FROM the city table GROUP BY AUTOMATIC HAVING COUNT ( last_update ) less than 158'
$ws.Range('B40').Value2 = 'IN Salaries table RETRIEVE Employees . ToDate and the count of ( FromDate ) AND THEN RETRIEVE Gender and Employees . BirthDate , the count of ( EmployeeNumber ) FROM Employees as Es where HireDate is equal to foo AND THEN natural join Salaries with Employees AND THEN GROUP BY AUTOMATIC'
$ws.Range('B41').Value2 = 'IN city where District equals foo DISPLAY the count of ( ID ) and country . Population AND THEN WHAT IS country . LocalName , city . LifeExpectancy and country . Capital FROM table country as cy where SurfaceArea less than 147 and Region is equal to bar AND THEN natural join city with country AND THEN GROUP BY AUTOMATIC'
$ws.Range('B42').Value2 = 'WHAT ARE facultyName and average ( capacity ) IN the courseoffering table AND THEN IN the term table PRESENT the sum of ( year ) and term . termPeriod and sum ( year ) AND THEN join courseoffering with term on courseoffering . termId equals term . id AND THEN GROUP BY AUTOMATIC'
$ws.Range('B43').Value2 = 'natural join DepartmentManager with departments AND THEN WHAT IS THE departments . DepartmentNumber , EmployeeNumber where DepartmentNumber is equal to bar FROM table DepartmentManager as Dr AND THEN IN table departments as ds WHAT IS count ( DepartmentName ) and DepartmentManager . DepartmentNumber , DepartmentName AND THEN GROUP BY AUTOMATIC'
$ws.Range('B44').Value2 = 'IN table rental where customer_id IN ( FIND rental_date FROM rental ) RETRIEVE the count of ( return_date ) , inventory_id AND THEN WHAT IS first_name and address_id IN the staff table AND THEN WHAT IS THE count ( address_id ) and rental . manager_staff_id where store_id equals bar and last_update IN ( FIND address_id FROM rental ) IN table store AND THEN natural join rental with staff AND THEN natural join staff with store AND THEN join store with staff on store . manager_staff_id equals staff . staff_id AND THEN GROUP BY AUTOMATIC'
$ws.Range('B45').Value2 = 'FROM city as cy where CountryCode IN ( SELECT city_id IN TABLE city ) DISPLAY city . CountryCode , count ( CountryCode ) , average ( Population ) , District HAVING SUM ( District ) equal to 48 GROUP BY AUTOMATIC'
$ws.Range('B46').Value2 = 'FROM table term as tm SHOW ME count ( id ) and termPeriod HAVING AVG ( termPeriod ) equal to 92 GROUP BY AUTOMATICALLY'
$ws.Range('B47').Value2 = 'IN DepartmentManager WHAT IS departments . DepartmentNumber , ToDate AND THEN WHAT ARE DepartmentManager . DepartmentNumber and the count of ( DepartmentName ) and count ( DepartmentName ) FROM departments AND THEN natural join DepartmentManager with departments AND THEN GROUP AUTOMATICALLY'
$ws.Range('B48').Value2 = 'natural join city with country AND THEN PRESENT country . country_id , ID and Population where country_id is equal to bar or city equals foo FROM city as cy AND THEN IN the country table WHAT ARE THE sum ( LifeExpectancy ) and LocalName and Code2 AND THEN HAVING AVG ( LocalName ) equal to 53 GROUP AUTOMATICALLY'
$ws.Range('B49').Value2 = 'natural join DepartmentManager with departments AND THEN FROM DepartmentManager as Dr WHAT IS THE ToDate , EmployeeNumber AND THEN FROM table departments as ds where DepartmentNumber is equal to bar and DepartmentNumber is equal to bar WHAT IS THE departments . DepartmentNumber , count ( DepartmentNumber ) AND THEN GROUP AUTOMATICALLY'
$ws.Range('B50').Value2 = 'SELECT count ( inventory_id ) and film_id IN the inventory table as iy AND THEN WHAT IS THE the count of ( film_id ) and length FROM the film table AND THEN natural join inventory with film AND THEN GROUP BY AUTOMATICALLY LIMIT 42'
$ws.Range('B51').Value2 = 'IN city table PRESENT sum ( Population ) , city . ID AND THEN SHOW ME city . Language and CountryCode and average ( Percentage ) , the sum of ( Percentage ) IN the countrylanguage table as ce where CountryCode is equal to foo and Language equals foo AND THEN natural join city with countrylanguage AND THEN HAVING COUNT ( Language ) greater than 43 GROUP AUTOMATICALLY'
$ws.Range('B52').Value2 = 'IN table courseoffering FIND term . startTime and term . courseId where roomId is equal to foo and capacity less than 490 AND THEN IN term DISPLAY the sum of ( year ) and term . termPeriod AND THEN join courseoffering with term on courseoffering . termId equals term . id AND THEN natural join term with term AND THEN natural join term with term AND THEN HAVING AVG ( termPeriod ) greater than 177 GROUP BY AUTOMATIC'
$ws.Range('B53').Value2 = 'natural join Titles with Employees AND THEN IN table Titles where FromDate IN ( SELECT Title FROM Titles ) GET Title , count ( ToDate ) AND THEN IN the Employees table as Es PRESENT count ( HireDate ) , LastName where Gender is equal to foo AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B54').Value2 = 'WHAT ARE the sum of ( Salary ) , the count of ( EmployeeNumber ) where Salary equals 1150 and FromDate equals foo IN the Salaries table as SsGROUP AUTOMATICALLY'
$ws.Range('B55').Value2 = 'natural join DepartmentManager with departments AND THEN IN the DepartmentManager table DISPLAY ToDate , departments . FromDate AND THEN RETRIEVE count ( DepartmentNumber ) , DepartmentManager . DepartmentNumber FROM ( WHAT ARE Salaries . EmployeeNumber , Salaries . Salary , the average of ( Salary ) FROM the Salaries table HAVING AVG ( EmployeeNumber ) greater than 31 GROUP BY AUTOMATICALLY ) as sq_alias AND THEN HAVING SUM ( DepartmentNumber ) equal to 67 GROUP BY AUTOMATIC'
$ws.Range('B56').Value2 = 'WHAT ARE THE DepartmentManager . DepartmentNumber , FromDate where FromDate is equal to bar and ToDate equals foo IN table DepartmentManager HAVING COUNT ( DepartmentNumber ) less than 108'
$ws.Range('B57').Value2 = 'natural join DepartmentManager with departments AND THEN FROM ( FROM table Titles RETRIEVE Titles . Title , ToDate , count ( EmployeeNumber ) GROUP BY AUTOMATIC ) as sq_alias where DepartmentNumber is equal to foo PRESENT count ( FromDate ) , departments . ToDate AND THEN PRESENT departments . DepartmentName , DepartmentNumber IN table departments AND THEN LIMIT 40 GROUP BY AUTOMATIC HAVING AVG ( DepartmentNumber ) equal to 199'
$ws.Range('B58').Value2 = 'natural join city with country AND THEN SHOW ME the sum of ( Population ) , country . ID where city_id IN ( DISPLAY city_id FROM TABLE country ) or ID is equal to foo IN city table AND THEN SHOW ME country and city . Region IN the country table as cy where country IN ( DISPLAY last_update FROM TABLE country ) AND THEN GROUP AUTOMATICALLY ORDER BY Region'
$ws.Range('B59').Value2 = 'WHAT IS city . Name and city . CountryCode , average ( Population ) and city . country_id where ID equals foo IN city GROUP BY AUTOMATICALLY'
$ws.Range('B60').Value2 = 'WHAT ARE THE country . city , the average of ( Population ) , country . ID IN table city as cy AND THEN DISPLAY country . IndepYear , country . GNPOld , city . HeadOfState and country . Capital FROM the country table AND THEN natural join city with country AND THEN GROUP AUTOMATIC'
$ws.Range('B61').Value2 = 'RETRIEVE the count of ( ToDate ) , ToDate where FromDate equals bar and EmployeeNumber is equal to bar IN DepartmentManager GROUP BY AUTOMATICALLY'
$ws.Range('B62').Value2 = 'FROM table city as cy SHOW ME city . District , the average of ( Population ) GROUP BY AUTOMATICALLY'
$ws.Range('B63').Value2 = 'natural join city with countrylanguage AND THEN WHAT IS THE the count of ( last_update ) , city_id where city is equal to bar or last_update is equal to bar IN city AND THEN FROM countrylanguage table SHOW ME IsOfficial and average ( Percentage ) and the count of ( Percentage ) , the average of ( Percentage ) where Percentage is equal to 408 or IsOfficial equals bar AND THEN GROUP AUTOMATIC LIMIT 33'
$ws.Range('B64').Value2 = 'natural join DepartmentEmployee with departments AND THEN IN DepartmentEmployee table where EmployeeNumber is equal to foo WHAT IS count ( ToDate ) , DepartmentEmployee . DepartmentNumber AND THEN SHOW ME DepartmentName and count ( DepartmentName ) FROM table departments AND THEN GROUP AUTOMATICALLY'
$ws.Range('B65').Value2 = 'IN ( IN courseoffering as cg where startTime is equal to bar or onDays is equal to foo WHAT ARE THE startTime and count ( capacity ) GROUP BY AUTOMATICALLY ) as sq_alias GET courseoffering . onDays and roomId , startTime , the average of ( capacity ) AND THEN WHAT IS THE course . deptId , id FROM course as ce AND THEN FROM the department table as dt GET id , course . departmentName AND THEN join courseoffering with course on courseoffering . courseId equals course . id AND THEN join course with department on course . deptId equals department . id AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B66').Value2 = 'natural join city with countrylanguage AND THEN FROM the city table WHAT IS THE the sum of ( Population ) , city , city_id AND THEN RETRIEVE city . IsOfficial , countrylanguage . Percentage and CountryCode FROM the countrylanguage table where IsOfficial is equal to bar AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B67').Value2 = 'DISPLAY ToDate , EmployeeNumber IN Titles HAVING SUM ( ToDate ) less than 64'
$ws.Range('B68').Value2 = 'FROM table DepartmentManager as Dr SELECT the count of ( DepartmentNumber ) , the count of ( DepartmentNumber ) and FromDate where DepartmentNumber is equal to foo or ToDate equals foo AND THEN IN the departments table WHAT ARE THE departments . DepartmentNumber and DepartmentManager . DepartmentName AND THEN natural join DepartmentManager with departments AND THEN GROUP BY AUTOMATIC'
$ws.Range('B69').Value2 = 'WHAT IS THE DepartmentManager . FromDate , the count of ( ToDate ) and DepartmentManager . ToDate IN the DepartmentManager table where EmployeeNumber IN ( SELECT DepartmentNumber FROM DepartmentManager ) or DepartmentNumber is equal to bar HAVING SUM ( FromDate ) less than 97 GROUP AUTOMATIC'
$ws.Range('B70').Value2 = 'IN courseoffering where termId is equal to bar SHOW ME the average of ( capacity ) and count ( termId ) GROUP AUTOMATICALLY'
$ws.Range('B71').Value2 = 'natural join DepartmentManager with departments AND THEN FROM DepartmentManager GET DepartmentManager . FromDate , count ( ToDate ) AND THEN IN departments as ds GET DepartmentName and DepartmentNumber where DepartmentName is equal to bar AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B72').Value2 = 'WHAT ARE THE Titles . FromDate , the count of ( EmployeeNumber ) where ToDate equals foo and Title equals foo FROM Titles GROUP BY AUTOMATICALLY'
$ws.Range('B73').Value2 = 'natural join address with city AND THEN natural join city with country AND THEN RETRIEVE city . last_update , postal_code IN ( SHOW ME last_update , payment_date IN table payment where customer_id equals foo ) as sq_alias AND THEN FROM table city FIND country . country_id and count ( country_id ) , city . District , city . Name where city_id is equal to foo or city_id IN ( DISPLAY last_update FROM TABLE city ) AND THEN DISPLAY count ( Code ) , IndepYear and city . GovernmentForm , address . Region IN country table AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B74').Value2 = 'FROM city table WHAT ARE THE Name and city . Population where Population greater than 1874 or last_update is equal to bar HAVING AVG ( Population ) greater than 133'
$ws.Range('B75').Value2 = 'FROM film where release_year equals foo WHAT ARE the count of ( length ) and title , rating GROUP BY AUTOMATICALLY HAVING SUM ( rating ) greater than 65'
$ws.Range('B76').Value2 = 'IN the city table as cy PRESENT country . ID and sum ( Population ) , sum ( Population ) AND THEN IN country table WHAT ARE THE average ( SurfaceArea ) , sum ( Capital ) and country . LocalName AND THEN natural join city with country AND THEN GROUP AUTOMATICALLY'
$ws.Range('B77').Value2 = 'DISPLAY courseoffering . id , capacity , sum ( capacity ) , sum ( capacity ) IN courseoffering as cg GROUP BY AUTOMATICALLY'
$ws.Range('B78').Value2 = 'PRESENT customer_id , rental_id , rental . return_date FROM rental as rl HAVING SUM ( return_date ) less than 85'
$ws.Range('B79').Value2 = 'natural join Titles with Employees AND THEN FIND count ( FromDate ) and Titles . Title where ToDate equals bar or EmployeeNumber IN ( FIND ToDate FROM Employees ) IN ( IN the DepartmentManager table as Dr where DepartmentNumber is equal to bar WHAT ARE the count of ( EmployeeNumber ) and the count of ( ToDate ) GROUP BY AUTOMATICALLY ) as sq_alias AND THEN IN the Employees table RETRIEVE FirstName , Employees . HireDate AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B80').Value2 = 'IN table Salaries FIND the sum of ( Salary ) , Salary where ToDate equals bar or ToDate IN ( FIND FromDate FROM Employees ) AND THEN FROM the Employees table as Es PRESENT Salaries . FirstName and Employees . HireDate AND THEN natural join Salaries with Employees AND THEN GROUP BY AUTOMATIC ORDER BY FirstName ASC'
$ws.Range('B81').Value2 = 'FROM ( FROM store where manager_staff_id equals bar SELECT the count of ( manager_staff_id ) , address_id and last_updateGROUP AUTOMATIC HAVING COUNT ( last_update ) equal to 35 ) as sq_alias where manager_staff_id equals bar or address_id equals foo WHAT IS THE count ( manager_staff_id ) , the count of ( store_id ) , last_update GROUP AUTOMATIC'
$ws.Range('B82').Value2 = 'WHAT ARE THE wheelchairSpaces and buildingId where id equals bar and capacity equals 1930 IN room HAVING COUNT ( wheelchairSpaces ) equal to 6 LIMIT 17'
$ws.Range('B83').Value2 = 'natural join city with countrylanguage AND THEN PRESENT count ( city ) , city . Population and average ( Population ) FROM table city where District is equal to foo AND THEN WHAT IS THE IsOfficial , CountryCode and the average of ( Percentage ) IN the countrylanguage table as ce where IsOfficial equals foo and Percentage less than 619 AND THEN GROUP AUTOMATICALLY'
$ws.Range('B84').Value2 = 'IN city as cy where Population less than 1751 and Population equals 1263 GET country_id and the count of ( Population ) GROUP BY AUTOMATICALLY'
$ws.Range('B85').Value2 = 'FROM table city as cy PRESENT countrylanguage . last_update and countrylanguage . District and sum ( Population ) AND THEN PRESENT Language and sum ( Percentage ) IN countrylanguage AND THEN natural join city with countrylanguage AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B86').Value2 = 'WHAT ARE count ( Title ) , ToDate FROM table Titles as TsGROUP BY AUTOMATICALLY'
$ws.Range('B87').Value2 = 'FROM city table WHAT IS Population and city . Name , last_update where ID IN ( SELECT city_id IN TABLE city ) or District equals foo LIMIT 17 GROUP AUTOMATIC HAVING AVG ( last_update ) equal to 136'
$ws.Range('B88').Value2 = 'FROM DepartmentManager table WHAT IS THE FromDate , DepartmentManager . ToDate and DepartmentManager . EmployeeNumber and count ( DepartmentNumber ) AND THEN PRESENT count ( DepartmentName ) , DepartmentNumber FROM the departments table AND THEN natural join DepartmentManager with departments AND THEN GROUP AUTOMATICALLY HAVING COUNT ( DepartmentNumber ) equal to 180'
$ws.Range('B89').Value2 = 'natural join city with country AND THEN FROM table city as cy WHAT ARE count ( CountryCode ) and country_id and country . last_update where Name equals foo and Name is equal to foo AND THEN FROM country table WHAT ARE sum ( Population ) and country . Population where LocalName equals foo and Code2 is equal to bar AND THEN GROUP AUTOMATIC HAVING AVG ( Population ) greater than 117'
$ws.Range('B90').Value2 = 'natural join Salaries with Employees AND THEN FROM Salaries WHAT IS Salaries . Salary , FromDate AND THEN PRESENT Employees . FirstName and EmployeeNumber FROM ( WHAT ARE DepartmentManager . DepartmentNumber , ToDate FROM table DepartmentManager where ToDate IN ( SELECT EmployeeNumber FROM TABLE DepartmentManager ) and ToDate IN ( DISPLAY EmployeeNumber FROM DepartmentManager )GROUP AUTOMATIC LIMIT 23 HAVING AVG ( ToDate ) equal to 58 ) as sq_alias where EmployeeNumber equals foo and EmployeeNumber equals bar AND THEN GROUP BY AUTOMATICALLY'
$ws.Range('B91').Value2 = 'SHOW ME staff . store_id and staff . username and staff . email and staff . last_update where username IN ( FIND last_update IN TABLE staff ) and address_id is equal to bar IN staff as sf LIMIT 48 HAVING AVG ( username ) greater than 2 GROUP AUTOMATIC'
$ws.Range('B92').Value2 = 'WHAT IS courseoffering . id , courseId and average ( capacity ) FROM courseoffering as cg GROUP BY AUTOMATIC HAVING SUM ( id ) greater than 199'
$ws.Range('B93').Value2 = 'RETRIEVE Salary , EmployeeNumber , sum ( Salary ) IN table Salaries GROUP BY AUTOMATIC'
$ws.Range('B94').Value2 = 'SHOW ME the count of ( EmployeeNumber ) , Employees . Title where Title is equal to foo and EmployeeNumber IN ( SELECT Title FROM TABLE Titles ) FROM Titles as Ts AND THEN GET count ( Gender ) , Employees . Gender , count ( HireDate ) IN Employees AND THEN natural join Titles with Employees AND THEN GROUP BY AUTOMATIC ORDER BY Gender ASC'
$ws.Range('B95').Value2 = 'IN the city table as cy RETRIEVE the average of ( Population ) , the average of ( Population ) where last_update equals foo and last_update IN ( FIND country_id IN TABLE city ) AND THEN GET countrylanguage . IsOfficial and average ( Percentage ) FROM countrylanguage as ce AND THEN natural join city with countrylanguage AND THEN GROUP BY AUTOMATIC'
$ws.Range('B96').Value2 = 'join course with department on course . deptId equals department . id AND THEN WHAT IS course . id and count ( units ) and units IN table course AND THEN WHAT IS department . id and the count of ( departmentName ) where departmentName is equal to bar FROM the department table AND THEN GROUP BY AUTOMATIC'
$ws.Range('B97').Value2 = 'FROM courseoffering as cg where id equals foo WHAT ARE THE endTime and the average of ( capacity ) LIMIT 10 GROUP BY AUTOMATIC'
$ws.Range('B98').Value2 = 'WHAT ARE THE sum ( capacity ) , the average of ( capacity ) and the count of ( facultyName ) and term . endTime where termId equals foo FROM the courseoffering table as cg AND THEN IN term table DISPLAY termPeriod , id and term . year , endDate where endDate is equal to foo AND THEN join courseoffering with term on courseoffering . termId equals term . id AND THEN natural join term with term AND THEN natural join term with term AND THEN GROUP AUTOMATICALLY'
$ws.Range('B99').Value2 = 'natural join DepartmentManager with departments AND THEN FROM DepartmentManager table where DepartmentNumber IN ( SELECT EmployeeNumber FROM TABLE departments ) and DepartmentNumber equals bar WHAT ARE THE DepartmentNumber , count ( DepartmentNumber ) AND THEN FROM departments table SHOW ME DepartmentManager . DepartmentName , departments . DepartmentNumber where DepartmentNumber equals foo AND THEN GROUP BY AUTOMATIC'
$ws.Range('B100').Value2 = 'natural join rental with staff AND THEN natural join staff with store AND THEN join store with staff on store . manager_staff_id equals staff . staff_id AND THEN PRESENT customer_id and count ( rental_id ) where inventory_id IN ( FIND return_date FROM TABLE rental ) or return_date is equal to foo IN table rental as rl AND THEN FROM the staff table as sf where username equals bar and last_name IN ( FIND address_id IN TABLE store ) WHAT ARE THE the average of ( active ) , last_name AND THEN RETRIEVE last_update , count ( address_id ) where address_id is equal to foo or last_update IN ( SELECT last_update FROM TABLE staff ) IN table store AND THEN GROUP AUTOMATIC HAVING COUNT ( last_update ) equal to 149'
$ws.Range('B101').Value2 = 'FROM city as cy FIND count ( country_id ) and count ( city ) GROUP BY AUTOMATICALLY'
